$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0372839506172839
$ws.Range("C2").Value = 0.000246913580246914
$ws.Range("D2").Value = 0.948641975308642
$ws.Range("F2").Value = 0.00148148148148148
$ws.Range("H2").Value = 0.000740740740740741
$ws.Range("I2").Value = 0.00148148148148148
$ws.Range("J2").Value = 0.947901234567901
$ws.Range("K2").Value = 0.000246913580246914
$ws.Range("L2").Value = 0.000493827160493827
$ws.Range("M2").Value = 0.907901234567901
$ws.Range("N2").Value = 0.000740740740740741
$ws.Range("P2").Value = 0.955061728395062
$ws.Range("R2").Value = 0.998518518518518
$ws.Range("T2").Value = 0.999259259259259
$ws.Range("U2").Value = 0.995061728395062
$ws.Range("V2").Value = 0.946913580246914
$ws.Range("W2").Value = 0.00271604938271605
$ws.Range("X2").Value = 0.0017283950617284
$ws.Range("C3").Value = 0.999259259259259
$ws.Range("D3").Value = 0.00567901234567901
$ws.Range("E3").Value = 0.0017283950617284
$ws.Range("G3").Value = 0.000987654320987654
$ws.Range("H3").Value = 0.000246913580246914
$ws.Range("J3").Value = 0.00839506172839506
$ws.Range("K3").Value = 0.997530864197531
$ws.Range("L3").Value = 0.998518518518518
$ws.Range("M3").Value = 0.00222222222222222
$ws.Range("O3").Value = 0.999753086419753
$ws.Range("P3").Value = 0.00641975308641975
$ws.Range("Q3").Value = 0.999259259259259
$ws.Range("R3").Value = 0.000740740740740741
$ws.Range("S3").Value = 0.909876543209877
$ws.Range("T3").Value = 0.000246913580246914
$ws.Range("U3").Value = 0.00345679012345679
$ws.Range("V3").Value = 0.000246913580246914
$ws.Range("W3").Value = 0.000493827160493827
$ws.Range("B4").Value = 0.942962962962963
$ws.Range("C4").Value = 0.000493827160493827
$ws.Range("D4").Value = 0.0133333333333333
$ws.Range("E4").Value = 0.000246913580246914
$ws.Range("F4").Value = 0.997530864197531
$ws.Range("G4").Value = 0.999012345679012
$ws.Range("H4").Value = 0.997777777777778
$ws.Range("I4").Value = 0.000246913580246914
$ws.Range("J4").Value = 0.0348148148148148
$ws.Range("K4").Value = 0.000987654320987654
$ws.Range("L4").Value = 0.000246913580246914
$ws.Range("M4").Value = 0.051358024691358
$ws.Range("N4").Value = 0.999012345679012
$ws.Range("P4").Value = 0.0365432098765432
$ws.Range("U4").Value = 0.000740740740740741
$ws.Range("V4").Value = 0.0461728395061728
$ws.Range("W4").Value = 0.996296296296296
$ws.Range("X4").Value = 0.997037037037037
$ws.Range("B5").Value = 0.0160493827160494
$ws.Range("D5").Value = 0.0276543209876543
$ws.Range("E5").Value = 0.998024691358025
$ws.Range("F5").Value = 0.000987654320987654
$ws.Range("H5").Value = 0.000246913580246914
$ws.Range("I5").Value = 0.998024691358025
$ws.Range("J5").Value = 0.00222222222222222
$ws.Range("K5").Value = 0.000987654320987654
$ws.Range("L5").Value = 0.000246913580246914
$ws.Range("M5").Value = 0.031358024691358
$ws.Range("N5").Value = 0.000246913580246914
$ws.Range("O5").Value = 0.000246913580246914
$ws.Range("Q5").Value = 0.000493827160493827
$ws.Range("R5").Value = 0.000246913580246914
$ws.Range("S5").Value = 0.0785185185185185
$ws.Range("T5").Value = 0.000246913580246914
$ws.Range("U5").Value = 0.000246913580246914
$ws.Range("W5").Value = 0.000246913580246914
$ws.Range("X5").Value = 0.000493827160493827
